# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, as published by the latest scrape.

$wb = $excel.ActiveWorkbook

# Row -> new F value (same updates apply identically on both sheets).
$updates = @{
    4  = 74
    5  = 12844
    6  = 68
    12 = 13706
    13 = 14161
    18 = 28
    22 = 1076
    23 = 110
    25 = 934
    26 = 5258
    28 = 281
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
